# Apply updated cryptocurrency data (prices and 1h volume %) per Fri Jul 26 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.231.03"
Set-TextValue "E2" "  +4.64%  "
Set-TextValue "D3" "3.239.53"
Set-TextValue "E3" "  +1.95%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "576.33"
Set-TextValue "E5" "  +2.04%  "
Set-TextValue "D6" "178.53"
Set-TextValue "E6" "  +5.63%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "E8" "  -1.26%  "
Set-TextValue "D9" "3.235.69"
Set-TextValue "E9" "  +1.99%  "
Set-TextValue "E10" "  +4.18%  "
Set-TextValue "E11" "  +1.53%  "
Set-TextValue "E12" "  +4.15%  "
Set-TextValue "D13" "3.804.32"
Set-TextValue "E13" "  +2.16%  "
Set-TextValue "E14" "  +0.74%  "
Set-TextValue "D15" "27.85"
Set-TextValue "E15" "  +1.83%  "
Set-TextValue "D16" "67.141.58"
Set-TextValue "E16" "  +4.52%  "
Set-TextValue "E17" "  +2.79%  "
Set-TextValue "D18" "3.242.47"
Set-TextValue "E18" "  +2.17%  "
Set-TextValue "D19" "5.85"
Set-TextValue "E19" "  +1.90%  "
Set-TextValue "D20" "13.28"
Set-TextValue "E20" "  +2.42%  "
Set-TextValue "D21" "372.96"
Set-TextValue "E21" "  +5.63%  "
Set-TextValue "E22" "  +5.69%  "
Set-TextValue "E23" "  +0.07%  "
Set-TextValue "D24" "70.99"
Set-TextValue "E24" "  +3.42%  "
Set-TextValue "E25" "  +0.69%  "
Set-TextValue "D26" "3.381.38"
Set-TextValue "E26" "  +2.24%  "
Set-TextValue "E27" "  -0.92%  "
Set-TextValue "D28" "9.90"
Set-TextValue "E28" "  +3.28%  "
Set-TextValue "D29" "0.180"
Set-TextValue "E29" "  +1.90%  "
Set-TextValue "E30" "  +0.45%  "
Set-TextValue "E31" "  +3.73%  "
Set-TextValue "D32" "5.62"
Set-TextValue "E32" "  +1.82%  "
Set-TextValue "D33" "22.54"
Set-TextValue "E33" "  +2.43%  "
Set-TextValue "B34" "Fetch.AI"
Set-TextValue "C34" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D34" "1.28"
Set-TextValue "E34" "  +6.39%  "
Set-TextValue "B35" "USDe"
Set-TextValue "C35" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  +0.03%  "
Set-TextValue "D36" "6.80"
Set-TextValue "E36" "  +2.40%  "
Set-TextValue "D37" "161.60"
Set-TextValue "E37" "  +5.03%  "
Set-TextValue "E38" "  +3.46%  "
Set-TextValue "D39" "0.855"
Set-TextValue "E39" "  +4.68%  "
Set-TextValue "E40" "  +9.86%  "
Set-TextValue "D41" "6.85"
Set-TextValue "E41" "  +14.47%  "
Set-TextValue "D42" "26.72"
Set-TextValue "E42" "  +2.77%  "
Set-TextValue "D43" "2.61"
Set-TextValue "E43" "  +4.72%  "
Set-TextValue "D44" "2.762.80"
Set-TextValue "E44" "  +5.67%  "
Set-TextValue "D45" "356.37"
Set-TextValue "E45" "  +11.42%  "
Set-TextValue "D46" "4.39"
Set-TextValue "E46" "  +5.06%  "
Set-TextValue "D47" "25.76"
Set-TextValue "E47" "  +8.77%  "
Set-TextValue "D48" "40.40"
Set-TextValue "E48" "  +2.60%  "
Set-TextValue "D49" "0.0672"
Set-TextValue "E49" "  +3.42%  "
Set-TextValue "D50" "0.0279"
Set-TextValue "E50" "  +3.21%  "
Set-TextValue "E51" "  +1.48%  "
